$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the trailing empty row (row 12), the R04 "password strength" test
# case (row 10, the old R04 becomes R05->R04 after this removal) and the
# L05 "Remember Me" test case (row 6). Deleting from the bottom up keeps
# the remaining row numbers stable while we work.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(6).Delete()

# Row height adjustments for the remaining rows.
$ws.Rows.Item(2).RowHeight = 13.55
$ws.Rows.Item(5).RowHeight = 31.75
$ws.Rows.Item(6).RowHeight = 35.6
$ws.Rows.Item(7).RowHeight = 40.5
$ws.Rows.Item(8).RowHeight = 44.3
$ws.Rows.Item(9).RowHeight = 62

# Column width adjustments (values chosen so the Excel character-width
# rounding lands on the desired stored width).
$ws.Columns.Item(3).ColumnWidth = 32.42857142857143
$ws.Columns.Item(4).ColumnWidth = 42
$ws.Columns.Item(5).ColumnWidth = 39.142857142857146

# Fix the "Page" column values that were wrong in the source data (the
# Register test cases were mistakenly labeled "Login").
$ws.Range("B6").Value = "Register"
$ws.Range("B7").Value = "Register"
$ws.Range("B8").Value = "Register"
$ws.Range("B9").Value = "Register"

# The old R04 (password strength) row was removed above, so the row that
# used to be R05 (confirm password) is renumbered to R04.
$ws.Range("A9").Value = "R04"
